$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "YM15 ZYU"
$ws.Range("B1").Value = "JAGUAR"
$ws.Range("C1").Value = "red"

$ws.Range("A1:C1").Select() | Out-Null
